$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 177; this shifts rows 177:285 down to 178:286
# (and therefore extends the used range to A1:T286, matching the target dimension).
$ws.Rows.Item(177).Insert()

# Populate the newly inserted row 177 with its data.
$ws.Range("A177").Value = 10
$ws.Range("B177").Value = "Vega Modelo de Temuco"
$ws.Range("C177").Value = "La Araucanía"
$ws.Range("D177").Value = 44582
$ws.Range("D177").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E177").Value = 9
$ws.Range("F177").Value = "Fruta"
$ws.Range("G177").Value = 100108
$ws.Range("H177").Value = "Tropicales y subtropicales"
$ws.Range("I177").Value = 100108002
$ws.Range("J177").Value = "Mango"
$ws.Range("K177").Value = "Sin especificar"
$ws.Range("L177").Value = "Primera"
$ws.Range("M177").Value = 155
$ws.Range("N177").Value = 8000
$ws.Range("O177").Value = 8000
$ws.Range("P177").Value = 8000
$ws.Range("Q177").Value = "$/bandeja 4 kilos"
$ws.Range("R177").Value = "Perú"
$ws.Range("S177").Value = 2000
$ws.Range("T177").Value = 4
